# Apply recalculated vm_pu results for the 380 kV case (slack bus voltage setpoint changed from 1.05 to 1.02 p.u.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
  2 = @(1.02, 1.011893977540996, 1.018668824758436, 0.9926147277508489, 1.010159802243678, 1.025148036254686, 1.017139824596926, 1.021514355129701, 0.9955398523336033, 1.013030835053877, 1.009859167459851)
  3 = @(1.02, 1.012742676172991, 1.019252157866555, 0.9936372048519304, 1.011653356366125, 1.025225366442524, 1.017622517048594, 1.021904161118772, 0.9963617723202692, 1.014326470288819, 1.010019660046266)
  4 = @(1.02, 1.013291873734151, 1.019629441121864, 0.9942998659930995, 1.012619753616572, 1.025273937752916, 1.017934273736274, 1.022155522651951, 0.9968940712668345, 1.015164318918161, 1.010123282965547)
  5 = @(1.02, 1.013522763807725, 1.019788008707978, 0.9945786998346017, 1.013026023867652, 1.025294005665575, 1.018065197691703, 1.02226098642184, 0.997117960005301, 1.015516431242016, 1.010166791679243)
  6 = @(1.02, 1.013561531689455, 1.019814630392504, 0.9946255319796338, 1.013094238377434, 1.025297354535035, 1.018087172266799, 1.022278682002502, 0.9971555583673453, 1.015575545568992, 1.010174093792023)
  7 = @(1.02, 1.013294958870013, 1.019631560077241, 0.9943035907982488, 1.012625182226292, 1.025274207282855, 1.017936023693364, 1.022156932684474, 0.9968970624462087, 1.015169024324061, 1.010123864544884)
  8 = @(1.02, 1.012180792168296, 1.018865999953535, 0.9929600610674301, 1.010664565241637, 1.025174473769727, 1.017303072203556, 1.021646271477642, 0.995817528259106, 1.013468809773028, 1.009913453503067)
  9 = @(1.02, 1.010217757752947, 1.017515711760702, 0.9906006454969559, 1.007209258081936, 1.024987515836221, 1.016183318709157, 1.020739790645921, 0.9939188001724441, 1.01046870692375, 1.009540953218389)
  10 = @(1.02, 1.008909271832444, 1.016614724042046, 0.989033133672735, 1.004905131989925, 1.024855359194961, 1.015433863775249, 1.020131044436919, 0.9926553831429383, 1.008465636153917, 1.0092914641762)
  11 = @(1.02, 1.008342735191567, 1.016224410062187, 0.988355674866747, 1.003907216480461, 1.024796355209481, 1.015108643503399, 1.019866407273273, 0.9921088820399291, 1.007597518317882, 1.009183159881577)
  12 = @(1.02, 1.008132305448103, 1.016079403925064, 0.9881042295826724, 1.003536508538471, 1.024774171507054, 1.014987736957245, 1.019767952547816, 0.9919059725120875, 1.007274939542145, 1.00914288977446)
  13 = @(1.02, 1.008177443013575, 1.016110509401744, 0.9881581567098651, 1.003616028445826, 1.024778942063459, 1.015013676583188, 1.019789078504398, 0.9919494934313052, 1.007344139343339, 1.009151529706384)
  14 = @(1.02, 1.008325340855406, 1.016212424337579, 0.9883348863814464, 1.003876574447558, 1.024794526942477, 1.015098651477161, 1.019858272168975, 0.9920921077337197, 1.007570856354658, 1.009179831979322)
  15 = @(1.02, 1.008416466556213, 1.016275214056761, 0.9884438009545853, 1.004037100276722, 1.024804093930987, 1.01515099336955, 1.01990088391865, 0.9921799884222134, 1.007710527991528, 1.009197264503954)
  16 = @(1.02, 1.008946871999846, 1.016640624176893, 0.9890781214508737, 1.004971355377168, 1.024859237635209, 1.015455432806856, 1.020148585534122, 0.9926916645766087, 1.008523233544716, 1.009298646211773)
  17 = @(1.02, 1.009279593900241, 1.016869788679458, 0.989476357848556, 1.005557327770815, 1.024893351686766, 1.015646211935092, 1.020303682483009, 0.9930127773699352, 1.00903281145397, 1.009362167006633)
  18 = @(1.02, 1.00947366941784, 1.017003438992152, 0.9897087662937556, 1.005899095310149, 1.024913078138886, 1.015757422496614, 1.020394046951582, 0.9932001317071769, 1.009329965171311, 1.009399191179697)
  19 = @(1.02, 1.0095398449041, 1.017049007303628, 0.9897880325774034, 1.006015625944765, 1.02491977521664, 1.015795330941595, 1.020424841762007, 0.9932640239640975, 1.009431274424335, 1.009411810980158)
  20 = @(1.02, 1.009243895525918, 1.016845203296628, 0.9894336180360679, 1.005494460612205, 1.024889709328952, 1.015625750139319, 1.020287052487381, 0.9929783193494215, 1.008978146301818, 1.009355354559869)
  21 = @(1.02, 1.008281788427032, 1.016182413643327, 0.9882828385668249, 1.003799851229344, 1.024789944951839, 1.015073631393566, 1.019837900681424, 0.9920501090198102, 1.007504097228222, 1.009171498798461)
  22 = @(1.02, 1.00767691521161, 1.015765540670806, 0.9875604150241495, 1.002734162211133, 1.024725674535363, 1.014725883492043, 1.019554594717462, 0.9914670000341481, 1.006576600730872, 1.009055663804114)
  23 = @(1.02, 1.007997565923, 1.015986546901795, 0.9879432794643023, 1.003299126926792, 1.024759891801555, 1.014910288918493, 1.01970486633495, 0.991776070289318, 1.007068352425448, 1.009117092620256)
  24 = @(1.02, 1.009260026070775, 1.016856312426878, 0.9894529299347244, 1.005522867636326, 1.024891355684189, 1.015634996153087, 1.020294567178502, 0.9929938892766442, 1.009002847359953, 1.009358432892972)
  25 = @(1.02, 1.010725214918578, 1.017864938254462, 0.9912096547607049, 1.008102620443033, 1.025037175658474, 1.016473323234465, 1.020974920064893, 0.9944092447426414, 1.011244816904892, 1.009637457646758)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $values[$i]
    }
}
